$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap values between row 2 and row 4 for columns D, J, K, L, M, O, P
$cols = @("D", "J", "K", "L", "M", "O", "P")

foreach ($col in $cols) {
    $cell2 = $ws.Range($col + "2")
    $cell4 = $ws.Range($col + "4")

    $val2 = $cell2.Value2
    $val4 = $cell4.Value2

    $cell2.Value2 = $val4
    $cell4.Value2 = $val2
}
